# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header - copy the "sum" header's formatting (bold, bordered, centered)
# onto the new "Save" header cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Save values for rows 2..46 (index 0 -> row 2, index 1 -> row 3, ...),
# derived from the source data.
$saveValues = @(
    1, 0, 0, 0, 0, 0, 0, 1, 0,
    1, 0, 1, 0, 0, 1, 0, 0, 0, 1,
    0, 0, 0, 1, 0, 0, 1, 0, 1, 0,
    0, 0, 1, 1, 0, 1, 0, 0, 0, 1,
    1, 0, 0, 1, 0, 0
)

for ($i = 0; $i -lt $saveValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
